$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the header merges to cover the two new week columns (H, I)
$ws.Range("F1:G1").UnMerge()
$ws.Range("F2:G2").UnMerge()
$ws.Range("F1:I1").Merge()
$ws.Range("F2:I2").Merge()

# Narrower week columns now that there are 4 of them instead of 2
# (COM ColumnWidth uses character-width units that don't map 1:1 onto the
# stored OOXML column width; 15.17 round-trips to the intended width of 16)
$ws.Columns.Item(6).ColumnWidth = 15.17
$ws.Columns.Item(7).ColumnWidth = 15.17

# New week header cells - copy the existing week-header formatting (fill/font/alignment)
# from G3 onto the new cells, then set their text
$ws.Range("G3").Copy($ws.Range("H3"))
$ws.Range("G3").Copy($ws.Range("I3"))
$ws.Range("H3").Value = "15/Jan - 21/Jan"
$ws.Range("I3").Value = "22/Jan - 28/Jan"

# Typo fix per commit diff
$ws.Range("C9").Value = "M2 - X@"
